# Weekly refresh of the "Fruta / hortaliza" price rows: column D (Fecha) and
# columns L..T (Calidad, Volumen, Precio mínimo/máximo/promedio, Unidad,
# Origen, Precio $/Kg, Kg/unidad) are rewritten per data row (rows 2-18,
# row 6 keeps its current weekly figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44650; L = "Primera"; M = 160; N = 31000; O = 32000; P = 31500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1575; T = 20 },
    @{ Row = 3;  D = 44650; L = "Segunda"; M = 250; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 4;  D = 44979; L = "Segunda"; M = 250; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 5;  D = 45021; L = "Segunda"; M = 250; N = 22000; O = 23000; P = 22500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1125; T = 20 },
    @{ Row = 7;  D = 45014; L = "Segunda"; M = 200; N = 24000; O = 25000; P = 24500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1225; T = 20 },
    @{ Row = 8;  D = 44671; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 9;  D = 44643; L = "Primera"; M = 160; N = 28000; O = 30000; P = 29000; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1450; T = 20 },
    @{ Row = 10; D = 45028; L = "Segunda"; M = 200; N = 21000; O = 22000; P = 21500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1075; T = 20 },
    @{ Row = 11; D = 44993; L = "Segunda"; M = 130; N = 25000; O = 26000; P = 25462; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1273; T = 20 },
    @{ Row = 12; D = 44972; L = "Segunda"; M = 140; N = 27000; O = 28000; P = 27429; Q = "$/caja 18 kilos"; R = "Región Metropolitana"; S = 1524; T = 18 },
    @{ Row = 13; D = 44679; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 14; D = 44679; L = "Tercera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1225; T = 20 },
    @{ Row = 15; D = 44636; L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 16; D = 44664; L = "Segunda"; M = 150; N = 29000; O = 30000; P = 29500; Q = "$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1639; T = 18 },
    @{ Row = 17; D = 44965; L = "Primera"; M = 100; N = 34000; O = 35000; P = 34600; Q = "$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1922; T = 18 },
    @{ Row = 18; D = 44965; L = "Segunda"; M = 120; N = 32000; O = 33000; P = 32333; Q = "$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1796; T = 18 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
}
